# Added Piranha Plant (no shooting mechanism)
#
# Appends a new "piranha plant" object block to the "Dynamic Object" sheet
# (rows 78-87), mirroring the existing Object/Sprite/Animation ID table
# layout used for the other dynamic objects (goomba, mushroom, superleaf, ...).

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("Dynamic Object")

# Row 78 - object header (Object name in column A) + first sub-state row
$ws.Range("A78").Value = "piranha plant"
$ws.Range("C78").Value = "top left closed"
$ws.Range("D78").Value = 130110
$ws.Range("E78").Value = "left move"
$ws.Range("F78").Value = 130000

# Row 79
$ws.Range("C79").Value = "top left opened"
$ws.Range("D79").Value = 130111
$ws.Range("E79").Value = "top left shoot"
$ws.Range("F79").Value = 130110

# Row 80
$ws.Range("C80").Value = "bottom left closed"
$ws.Range("D80").Value = 130210
$ws.Range("E80").Value = "bottom left shoot"
$ws.Range("F80").Value = 130210

# Row 81
$ws.Range("C81").Value = "bottom left opened"
$ws.Range("D81").Value = 130211

# Row 82
$ws.Range("C82").Value = "top right closed"
$ws.Range("D82").Value = 130120
$ws.Range("E82").Value = "right move"
$ws.Range("F82").Value = 130010

# Row 83
$ws.Range("C83").Value = "top right opened"
$ws.Range("D83").Value = 130121
$ws.Range("E83").Value = "top right shoot"
$ws.Range("F83").Value = 130120

# Row 84
$ws.Range("C84").Value = "bottom right closed"
$ws.Range("D84").Value = 130220
$ws.Range("E84").Value = "bottom right shoot"
$ws.Range("F84").Value = 130220

# Row 85
$ws.Range("C85").Value = "bottom right opened"
$ws.Range("D85").Value = 130221

# Row 86
$ws.Range("C86").Value = "upright closed"
$ws.Range("D86").Value = 130000

# Row 87
$ws.Range("C87").Value = "upright opened"
$ws.Range("D87").Value = 130001

# Bring the "Dynamic Object" sheet to the front / make it the active tab,
# and leave the selection sitting on the new block, matching the author's
# final view state.
$ws.Activate() | Out-Null
$ws.Range("E81").Select() | Out-Null
